$d = $word.ActiveDocument

# --- Update the first paragraph's formatting: add a paragraph border
#     (5-twip spacing on all four sides) and change the left indent
#     from 120 twips (6pt) to 225 twips (11.25pt). ---
$p1 = $d.Paragraphs.Item(1)
$pf1 = $p1.Range.ParagraphFormat
$borders1 = $pf1.Borders
$borders1.DistanceFromTop = 5
$borders1.DistanceFromLeft = 5
$borders1.DistanceFromBottom = 5
$borders1.DistanceFromRight = 5
$pf1.LeftIndent = 11.25

# --- Replace the bookmark ID text in that same paragraph. ---
$findRange = $d.Content
$old = "**ID__AFFARS_pgi_5309_topic_11__ID**"
$new = "**ID__AFFARS_SMC_PGI_5309_5__ID**"
$found = $findRange.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $findRange.Text = $new
}

# --- The paragraph used to end with a separate run containing just a
#     trailing space; drop that trailing space so the paragraph holds
#     only the single run with the new ID text. ---
$p1 = $d.Paragraphs.Item(1)
$p1Range = $p1.Range
$lastChar = $d.Range($p1Range.End - 2, $p1Range.End - 1)
if ($lastChar.Text -eq " ") {
    $lastChar.Text = ""
}
